# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-16 13:34:48
#
# The "Recorded By" column (G) lists the users who touched a session record.
# Two cosmetic ordering issues are fixed here:
#   1. "dnasr281@gmail.com, System"              -> "System, dnasr281@gmail.com"
#   2. "backup@backdoor.com, system, System"     -> "backup@backdoor.com, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "backup@backdoor.com, system, System") {
        $cell.Value = "backup@backdoor.com, System, system"
    }
}
